$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, pushing current rows 31-32 down to 32-33.
$ws.Rows.Item(31).Insert()

# Copy formatting (style) of the date cell from the row above into the new row's date cell.
$ws.Cells.Item(30, 4).Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4122) | Out-Null

# Populate the new row 31 with the new data entry.
$ws.Cells.Item(31, 1).Value = 9
$ws.Cells.Item(31, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(31, 3).Value = "Metropolitana"
$ws.Cells.Item(31, 4).Value = 44706
$ws.Cells.Item(31, 5).Value = 13
$ws.Cells.Item(31, 6).Value = 100112035
$ws.Cells.Item(31, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 30
$ws.Cells.Item(31, 11).Value = 26000
$ws.Cells.Item(31, 12).Value = 26000
$ws.Cells.Item(31, 13).Value = 26000
$ws.Cells.Item(31, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(31, 15).Value = "Hijuelas"
$ws.Cells.Item(31, 16).Value = 1733
$ws.Cells.Item(31, 17).Value = 15
$ws.Cells.Item(31, 18).Value = "Hortaliza"
